# Commit: "renamed some variables for clarity separated pro(pro-com creation
# method into two separate methods"
#
# On the trafo3w_std_types sheet this reorders the rated-power columns
# (sn_hv_mva/sn_mv_mva/sn_lv_mva) ahead of the rated-voltage columns
# (vn_hv_kv/vn_mv_kv/vn_lv_kv) and introduces a new vector_group column
# (value "YN0yn0yn0" for both existing standard-type rows), pushing the
# tap_* columns one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trafo3w_std_types")

# ---- Header row (row 1), columns B..W --------------------------------------
$headers = @(
    "sn_hv_mva",
    "sn_mv_mva",
    "sn_lv_mva",
    "vn_hv_kv",
    "vn_mv_kv",
    "vn_lv_kv",
    "vk_hv_percent",
    "vk_mv_percent",
    "vk_lv_percent",
    "vkr_hv_percent",
    "vkr_mv_percent",
    "vkr_lv_percent",
    "pfe_kw",
    "i0_percent",
    "shift_mv_degree",
    "shift_lv_degree",
    "vector_group",
    "tap_side",
    "tap_neutral",
    "tap_min",
    "tap_max",
    "tap_step_percent"
)

# New rightmost header cell (W1) needs the same bold/border/center format as
# its neighbors - clone it from V1 before the values shift the range out.
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# ---- Data rows ---------------------------------------------------------------
# Row 2: "63/25/38 MVA 110/20/10 kV"
$row2 = @(63, 25, 38, 110, 20, 10, 10.4, 10.4, 10.4, 0.28, 0.32, 0.35, 35, 0.89, 0, 0, "YN0yn0yn0", "hv", 0, -10, 10, 1.2)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $row2[$i]
}

# Row 3: "63/25/38 MVA 110/10/10 kV"
$row3 = @(63, 25, 38, 110, 10, 10, 10.4, 10.4, 10.4, 0.28, 0.32, 0.35, 35, 0.89, 0, 0, "YN0yn0yn0", "hv", 0, -10, 10, 1.2)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 2 + $i).Value = $row3[$i]
}

Write-Output "trafo3w_std_types updated"
